# Scheduled-runner refresh of market-price-derived columns (H:N) across all
# crafting-job sheets: currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ) and the
# resulting LeveProfit(NQ/HQ). Values below mirror the latest pull; a few
# rows also gain/lose an H:N cell entirely where a HQ/NQ price path that was
# previously zero (and therefore profit not computed/omitted) now resolves
# (or vice versa).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2453.0833
$ws.Range("I19").Value = 2184.4285
$ws.Range("J19").Value = 2829.2
$ws.Range("K19").Value = 2184.4285
$ws.Range("L19").Value = 2829.2
$ws.Range("M19").Value = -2009.4285
$ws.Range("N19").Value = -3179.2

$ws.Range("H132").Value = 2358.652
$ws.Range("I132").Value = 2028.2941
$ws.Range("K132").Value = 6084.8823
$ws.Range("M132").Value = -3554.8823

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3973.25
$ws.Range("I2").Value = 1950
$ws.Range("J2").Value = 5996.5
$ws.Range("K2").Value = 1950
$ws.Range("L2").Value = 5996.5
$ws.Range("M2").Value = -1837
$ws.Range("N2").Value = -6222.5

$ws.Range("H88").Value = 3555.889
$ws.Range("I88").Value = 3668
$ws.Range("K88").Value = 3668
$ws.Range("M88").Value = -3262

$ws.Range("H91").Value = 3555.889
$ws.Range("I91").Value = 3668
$ws.Range("K91").Value = 3668
$ws.Range("M91").Value = -2264

$ws.Range("H116").Value = 3973.25
$ws.Range("I116").Value = 1950
$ws.Range("J116").Value = 5996.5
$ws.Range("K116").Value = 1950
$ws.Range("L116").Value = 5996.5
$ws.Range("M116").Value = 344
$ws.Range("N116").Value = -10584.5

$ws.Range("H132").Value = 1550
$ws.Range("I132").Value = 1550
$ws.Range("K132").Value = 4650
$ws.Range("M132").Value = -2120

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3973.25
$ws.Range("I3").Value = 1950
$ws.Range("J3").Value = 5996.5
$ws.Range("K3").Value = 1950
$ws.Range("L3").Value = 5996.5
$ws.Range("M3").Value = -1836
$ws.Range("N3").Value = -6224.5

$ws.Range("H99").Value = 2127.1428
$ws.Range("I99").Value = 2031.75
$ws.Range("K99").Value = 2031.75
$ws.Range("M99").Value = -533.75

$ws.Range("H105").Value = 3114.1428
$ws.Range("I105").Value = 3114.1428
$ws.Range("K105").Value = 3114.1428
$ws.Range("M105").Value = -1367.1428

$ws.Range("H134").Value = 1066.1305
$ws.Range("I134").Value = 1100.9546
$ws.Range("K134").Value = 3302.8638
$ws.Range("M134").Value = -767.8638000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 5654.4546
$ws.Range("I99").Value = 5519.9
$ws.Range("J99").Value = 7000
$ws.Range("K99").Value = 5519.9
$ws.Range("L99").Value = 7000
$ws.Range("M99").Value = -4021.9
$ws.Range("N99").Value = -9996

$ws.Range("H126").Value = 5654.4546
$ws.Range("I126").Value = 5519.9
$ws.Range("J126").Value = 7000
$ws.Range("K126").Value = 16559.7
$ws.Range("L126").Value = 21000
$ws.Range("M126").Value = -14089.7
$ws.Range("N126").Value = -25940

$ws.Range("H132").Value = 2302.1667
$ws.Range("I132").Value = 2302.1667
$ws.Range("K132").Value = 6906.500100000001
$ws.Range("M132").Value = -4376.500100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H54").Value = 2000
$ws.Range("J54").Value = 2000
$ws.Range("L54").Value = 6000
$ws.Range("N54").Value = -7118

$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()

$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()

$ws.Range("H69").Value = 1097.5
$ws.Range("I69").Value = 518.2
$ws.Range("J69").Value = 3994
$ws.Range("K69").Value = 1554.6
$ws.Range("L69").Value = 11982
$ws.Range("M69").Value = -743.6000000000001
$ws.Range("N69").Value = -13604

$ws.Range("H72").Value = 1097.5
$ws.Range("I72").Value = 518.2
$ws.Range("J72").Value = 3994
$ws.Range("K72").Value = 4663.8
$ws.Range("L72").Value = 35946
$ws.Range("M72").Value = -607.8000000000002
$ws.Range("N72").Value = -44058

$ws.Range("H103").Value = 28733.223
$ws.Range("I103").Value = 50474.8
$ws.Range("J103").Value = 1556.25
$ws.Range("K103").Value = 151424.4
$ws.Range("L103").Value = 4668.75
$ws.Range("M103").Value = -150545.4
$ws.Range("N103").Value = -6426.75

$ws.Range("H131").Value = 1099.3334
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

$ws.Range("H140").Value = 517.5
$ws.Range("I140").Value = 517.5
$ws.Range("K140").Value = 1552.5
$ws.Range("M140").Value = 3627.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2480.6667
$ws.Range("I102").Value = 2342
$ws.Range("J102").Value = 2896.6667
$ws.Range("K102").Value = 2342
$ws.Range("L102").Value = 2896.6667
$ws.Range("M102").Value = -720
$ws.Range("N102").Value = -6140.6667

$ws.Range("H113").Value = 1389
$ws.Range("J113").Value = 1389
$ws.Range("L113").Value = 1389
$ws.Range("N113").Value = -5729

$ws.Range("H122").Value = 6249.75
$ws.Range("I122").Value = 5000
$ws.Range("J122").Value = 9999
$ws.Range("K122").Value = 15000
$ws.Range("L122").Value = 29997
$ws.Range("M122").Value = -12550
$ws.Range("N122").Value = -34897

$ws.Range("H123").Value = 30124.75
$ws.Range("J123").Value = 30124.75
$ws.Range("L123").Value = 30124.75
$ws.Range("N123").Value = -35024.75

$ws.Range("H132").Value = 2467.4285
$ws.Range("I132").Value = 2471.5386
$ws.Range("K132").Value = 7414.6158
$ws.Range("M132").Value = -4884.6158

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 738.3333
$ws.Range("I16").Value = 738.3333
$ws.Range("K16").Value = 738.3333
$ws.Range("M16").Value = -568.3333

$ws.Range("H22").Value = 1866.3334
$ws.Range("I22").Value = 499.5
$ws.Range("K22").Value = 499.5
$ws.Range("M22").Value = -204.5

$ws.Range("H27").Value = 1866.3334
$ws.Range("I27").Value = 499.5
$ws.Range("K27").Value = 499.5
$ws.Range("M27").Value = -392.5

$ws.Range("H46").Value = 4218.5806
$ws.Range("I46").Value = 3748.6875
$ws.Range("J46").Value = 4719.8
$ws.Range("K46").Value = 3748.6875
$ws.Range("L46").Value = 4719.8
$ws.Range("M46").Value = -3560.6875
$ws.Range("N46").Value = -5095.8

$ws.Range("H55").Value = 485.35294
$ws.Range("I55").Value = 45
$ws.Range("J55").Value = 544.06665
$ws.Range("K55").Value = 45
$ws.Range("L55").Value = 544.06665
$ws.Range("M55").Value = 128
$ws.Range("N55").Value = -890.06665

$ws.Range("H61").Value = 915
$ws.Range("I61").Value = 915
$ws.Range("K61").Value = 915
$ws.Range("M61").Value = -713

$ws.Range("H113").Value = 915
$ws.Range("I113").Value = 915
$ws.Range("K113").Value = 915
$ws.Range("M113").Value = 1255

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 4097.5

$ws.Range("H107").Value = 609.6
$ws.Range("I107").Value = 611.25
$ws.Range("K107").Value = 1833.75
$ws.Range("M107").Value = 86.25

$ws.Range("H136").Value = 1087.5714
$ws.Range("I136").Value = 1147.7693
$ws.Range("K136").Value = 3443.3079
$ws.Range("M136").Value = -893.3078999999998

$ws.Range("H140").Value = 79998.5
$ws.Range("I140").Value = 79997
$ws.Range("K140").Value = 79997
$ws.Range("M140").Value = -74817

$ws.Range("H141").Value = 89999.5
$ws.Range("I141").Value = 89999
$ws.Range("K141").Value = 89999
$ws.Range("M141").Value = -84819
